$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "UNO International Corp / UNO Premier Philippines International Corp"
$ws.Range("A2").Value = "TIN: 006-505-014 VAT / 008-902-716 VAT"

$ws.Range("A49").Select()
